# Convert the three M2Doc field-code tokens ("m: for ...", "m:v.name",
# "m:endfor") from real Word fields (fldChar begin/instrText/fldChar end)
# into plain literal text using the "{m: ...}" token syntax, as done by
# the TokenIteratorFieldRewriterSplit parser update.
#
# Each field is removed with Field.Delete() (which also strips the
# fldChar/instrText runs) and replaced by inserting the equivalent
# literal text directly into the paragraph. The hidden "_GoBack"
# bookmark that sat inside the first field's code is recreated at the
# same textual position once the literal text is in place.

$d = $word.ActiveDocument

# --- Field 1: " m: for v | self.eClassifiers " -> lives in paragraph 2 ---
# --- Field 2: " m:v.name "                     -> lives in paragraph 3 ---
# --- Field 3: " m:endfor "                     -> lives in paragraph 4 ---
# Deleting always removes whatever field is currently first, in document
# order, so the same call grabs each one in turn.
$f1 = $d.Fields.Item(1)
$f1.Delete()
$f2 = $d.Fields.Item(1)
$f2.Delete()
$f3 = $d.Fields.Item(1)
$f3.Delete()

# --- Paragraph 2: " " -> " {m: for v | self.eClassifiers}" -------------
# (the leading space run is untouched; the field content is rebuilt as
# plain text, with the _GoBack bookmark re-inserted right after "{m: ")
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start
$insertPoint2 = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$insertPoint2.InsertBefore("{m: for v | self.eClassifiers}")

$bookmarkPos = $p2Start + 5   # length of " {m: " before the bookmark
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- Paragraph 3: "name = ," -> "name = {m:v.name}," -------------------
$p3 = $d.Paragraphs.Item(3)
$insertPoint3 = $d.Range($p3.Range.End - 2, $p3.Range.End - 2)
$insertPoint3.InsertBefore("{m:v.name}")

# --- Paragraph 4: "" -> "{m:endfor}" ------------------------------------
$p4 = $d.Paragraphs.Item(4)
$insertPoint4 = $d.Range($p4.Range.End - 1, $p4.Range.End - 1)
$insertPoint4.InsertBefore("{m:endfor}")
